$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row above the current row 64 ("rdg/wit" row) ---
# This shifts the old rows 64..71 down to 65..72 (old row 71, a blank
# all-style-1 row, ends up at row 72 and gets removed further below).
$ws.Rows.Item(64).Insert()

# Copy the formatting (styles) of the standard "odd" data row (row 5,
# style pattern s=5,8,5,3) onto the four cells of the new row 64.
$ws.Range("A5").Copy()
$ws.Range("A64").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("B64").PasteSpecial(-4122)
$ws.Range("C5").Copy()
$ws.Range("C64").PasteSpecial(-4122)
$ws.Range("D5").Copy()
$ws.Range("D64").PasteSpecial(-4122)

# Fill in the new row's content: "type" / explanation of identical-missing markup
$ws.Range("C64").Value = "type"
$ws.Range("D64").Value = '"identical" gleicher Text in beiden Varianten; "missing" fehlt in der anderen Version'

# Match the taller row height used for this note row
$ws.Rows.Item(64).RowHeight = 32

# --- Remove the old trailing blank row (now shifted to row 72) ---
$ws.Rows.Item(72).Delete()

# --- Append a new blank row at the very end of the sheet data (row 161) ---
$ws.Range("A160:D160").Copy()
$ws.Range("A161:D161").PasteSpecial(-4122)

# --- Update the view state (scroll position & active selection) ---
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.ScrollRow = 55
$ws.Range("C65").Select()
